# Applies the scheduled-runner market-data refresh described in the commit.
# For each affected leve row, columns H:N (currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) are
# refreshed with newly-fetched market values. A handful of rows gain or lose a cell
# entirely (e.g. a profit column that did not previously apply to that leve).

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H15" = 288
    "I15" = 288
    "K15" = 864
    "M15" = -695
    "H19" = 201076.11
    "J19" = 821.5
    "L19" = 821.5
    "N19" = -1171.5
    "H33" = 261.8095
    "I33" = 151.58333
    "J33" = 408.77777
    "K33" = 151.58333
    "L33" = 408.77777
    "M33" = 77.41667000000001
    "N33" = -866.7777699999999
    "H116" = 3709.64
    "I116" = 3564.6875
    "J116" = 3967.3333
    "K116" = 3564.6875
    "L116" = 3967.3333
    "M116" = -122.6875
    "N116" = -10851.3333
    "H132" = 1591.2444
    "I132" = 1591.2444
    "J132" = 0
    "K132" = 4773.733200000001
    "L132" = 0
    "M132" = -2243.733200000001
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("N132")) {
    $ws.Range($ref).ClearContents()
}

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H18" = 80013
    "J18" = 80013
    "L18" = 80013
    "N18" = -80657
    "H61" = 1831.5
    "I61" = 1748.0358
    "J61" = 3000
    "K61" = 1748.0358
    "L61" = 3000
    "M61" = -1536.0358
    "N61" = -3424
    "H74" = 1039.6268
    "I74" = 961.4262
    "J74" = 1834.6666
    "K74" = 961.4262
    "L74" = 1834.6666
    "M74" = -87.42619999999999
    "N74" = -3582.6666
    "H77" = 1039.6268
    "I77" = 961.4262
    "J77" = 1834.6666
    "K77" = 4807.131
    "L77" = 9173.333000000001
    "M77" = -439.1310000000003
    "N77" = -17909.333
    "H122" = 3519.7144
    "I122" = 2656
    "J122" = 4671.3335
    "K122" = 7968
    "L122" = 14014.0005
    "M122" = -5518
    "N122" = -18914.0005
    "H132" = 1632.8541
    "I132" = 1427.921
    "J132" = 2411.6
    "K132" = 4283.763
    "L132" = 7234.799999999999
    "M132" = -1753.763
    "N132" = -12294.8
    "H136" = 1831.5
    "I136" = 1748.0358
    "J136" = 3000
    "K136" = 5244.107400000001
    "L136" = 9000
    "M136" = -2694.107400000001
    "N136" = -14100
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H134" = 24820.842
    "I134" = 1820.742
    "J134" = 79667.234
    "K134" = 5462.226
    "L134" = 239001.702
    "M134" = -2927.226
    "N134" = -244071.702
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H22" = 7240973.5
    "I22" = 21717172
    "J22" = 2875
    "K22" = 65151516
    "L22" = 8625
    "M22" = -65151347
    "N22" = -8963
    "H27" = 7240973.5
    "I27" = 21717172
    "J27" = 2875
    "K27" = 65151516
    "L27" = 8625
    "M27" = -65151414
    "N27" = -8829
    "H68" = 1479.4546
    "I68" = 433
    "J68" = 1871.875
    "K68" = 1299
    "L68" = 5615.625
    "M68" = -488
    "N68" = -7237.625
    "H71" = 1479.4546
    "I71" = 433
    "J71" = 1871.875
    "K71" = 3897
    "L71" = 16846.875
    "M71" = 159
    "N71" = -24958.875
    "H97" = 4464585.5
    "I97" = 5952705.5
    "J97" = 225
    "K97" = 17858116.5
    "L97" = 675
    "M97" = -17857620.5
    "N97" = -1667
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H17" = 80009
    "J17" = 80009
    "L17" = 80009
    "N17" = -80345
    "H97" = 761.0833
    "I97" = 888.46155
    "J97" = 429.9
    "K97" = 888.46155
    "L97" = 429.9
    "M97" = -392.46155
    "N97" = -1421.9
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H61" = 1211.1765
    "I61" = 1206.4286
    "K61" = 1206.4286
    "M61" = -1004.4286
    "H68" = 2056.6
    "I68" = 2079.5293
    "J68" = 1926.6666
    "K68" = 2079.5293
    "L68" = 1926.6666
    "M68" = -1330.5293
    "N68" = -3424.6666
    "H71" = 2056.6
    "I71" = 2079.5293
    "J71" = 1926.6666
    "K71" = 10397.6465
    "L71" = 9633.333000000001
    "M71" = -6653.646500000001
    "N71" = -17121.333
    "H82" = 1755.1923
    "I82" = 1744.381
    "J82" = 1800.6
    "K82" = 1744.381
    "L82" = 1800.6
    "M82" = -1383.381
    "N82" = -2522.6
    "H85" = 1755.1923
    "I85" = 1744.381
    "J85" = 1800.6
    "K85" = 1744.381
    "L85" = 1800.6
    "M85" = -496.3810000000001
    "N85" = -4296.6
    "H93" = 2553.7368
    "I93" = 2423.5652
    "J93" = 2753.3333
    "K93" = 2423.5652
    "L93" = 2753.3333
    "M93" = -1175.5652
    "N93" = -5249.3333
    "H113" = 1211.1765
    "I113" = 1206.4286
    "K113" = 1206.4286
    "M113" = 963.5714
    "H122" = 1740.3
    "I122" = 1629
    "J122" = 2000
    "K122" = 4887
    "L122" = 6000
    "M122" = -2437
    "N122" = -10900
    "H136" = 6226.5
    "I136" = 2258.8572
    "K136" = 6776.571599999999
    "M136" = -4226.571599999999
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H18" = 0
    "J18" = 0
    "L18" = 0
    "H107" = 1110.1154
    "I107" = 1025.6957
    "J107" = 1757.3334
    "K107" = 3077.0871
    "L107" = 5272.0002
    "M107" = -1157.0871
    "N107" = -9112.0002
    "H113" = 641.6667
    "I113" = 513.2
    "J113" = 898.6
    "K113" = 1539.6
    "L113" = 2695.8
    "M113" = 630.3999999999999
    "N113" = -7035.8
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("N18")) {
    $ws.Range($ref).ClearContents()
}

